$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.1477750351608889
$ws.Range("J2").Value = 0.1477750351608889
$ws.Range("M2").Value = 3.867218333333334
$ws.Range("N2").Value = 11.601655
$ws.Range("O2").Value = 0.1566152977872902
$ws.Range("P2").Value = 0.1566152977872902
$ws.Range("Q2").Value = 0.2819936936483333
$ws.Range("R2").Value = 2.537943242835
$ws.Range("S2").Value = 0.0231438311372499
$ws.Range("T2").Value = 0.0231438311372499

# Row 3
$ws.Range("I3").Value = 0.1477750351608889
$ws.Range("J3").Value = 0.1477750351608889
$ws.Range("N3").Value = 33.813685
$ws.Range("O3").Value = 0.4564642152831324
$ws.Range("P3").Value = 0.4564642152831324
$ws.Range("Q3").Value = 0.8218866988383333
$ws.Range("R3").Value = 7.396980289545001
$ws.Range("S3").Value = 0.06745401546315245
$ws.Range("T3").Value = 0.06745401546315245

# Row 4
$ws.Range("I4").Value = 0.1477750351608889
$ws.Range("J4").Value = 0.1477750351608889
$ws.Range("M4").Value = 5.654344666666667
$ws.Range("N4").Value = 16.963034
$ws.Range("O4").Value = 0.2289906587711778
$ws.Range("P4").Value = 0.2289906587711778
$ws.Range("Q4").Value = 0.4123091587486666
$ws.Range("R4").Value = 3.710782428738
$ws.Range("S4").Value = 0.03383910265142592
$ws.Range("T4").Value = 0.03383910265142592

# Row 5
$ws.Range("I5").Value = 0.1477750351608889
$ws.Range("J5").Value = 0.1477750351608889
$ws.Range("M5").Value = 0.819389
$ws.Range("N5").Value = 2.458167
$ws.Range("O5").Value = 0.03318376186120772
$ws.Range("P5").Value = 0.03318376186120772
$ws.Range("Q5").Value = 0.059749026491
$ws.Range("R5").Value = 0.537741238419
$ws.Range("S5").Value = 0.004903731575810535
$ws.Range("T5").Value = 0.004903731575810535

# Row 6
$ws.Range("I6").Value = 0.1477750351608889
$ws.Range("J6").Value = 0.1477750351608889
$ws.Range("M6").Value = 3.080288333333333
$ws.Range("N6").Value = 9.240864999999999
$ws.Range("O6").Value = 0.1247460662971919
$ws.Range("P6").Value = 0.1247460662971919
$ws.Range("Q6").Value = 0.2246115449783333
$ws.Range("R6").Value = 2.021503904805
$ws.Range("S6").Value = 0.01843435433325011
$ws.Range("T6").Value = 0.01843435433325011

# Row 7
$ws.Range("G7").Value = 0.420527
$ws.Range("H7").Value = 1.261581
$ws.Range("I7").Value = 0.852224964839111
$ws.Range("J7").Value = 0.852224964839111
$ws.Range("M7").Value = 3.867218333333334
$ws.Range("N7").Value = 11.601655
$ws.Range("O7").Value = 0.1566152977872902
$ws.Range("P7").Value = 0.1566152977872902
$ws.Range("Q7").Value = 1.626269724061667
$ws.Range("R7").Value = 14.636427516555
$ws.Range("S7").Value = 0.1334714666500403
$ws.Range("T7").Value = 0.1334714666500403

# Row 8
$ws.Range("G8").Value = 0.420527
$ws.Range("H8").Value = 1.261581
$ws.Range("I8").Value = 0.852224964839111
$ws.Range("J8").Value = 0.852224964839111
$ws.Range("N8").Value = 33.813685
$ws.Range("O8").Value = 0.4564642152831324
$ws.Range("P8").Value = 0.4564642152831324
$ws.Range("Q8").Value = 4.739855837331667
$ws.Range("R8").Value = 42.658702535985
$ws.Range("S8").Value = 0.3890101998199799
$ws.Range("T8").Value = 0.3890101998199799

# Row 9
$ws.Range("G9").Value = 0.420527
$ws.Range("H9").Value = 1.261581
$ws.Range("I9").Value = 0.852224964839111
$ws.Range("J9").Value = 0.852224964839111
$ws.Range("M9").Value = 5.654344666666667
$ws.Range("N9").Value = 16.963034
$ws.Range("O9").Value = 0.2289906587711778
$ws.Range("P9").Value = 0.2289906587711778
$ws.Range("Q9").Value = 2.377804599639334
$ws.Range("R9").Value = 21.400241396754
$ws.Range("S9").Value = 0.1951515561197519
$ws.Range("T9").Value = 0.1951515561197519

# Row 10
$ws.Range("G10").Value = 0.420527
$ws.Range("H10").Value = 1.261581
$ws.Range("I10").Value = 0.852224964839111
$ws.Range("J10").Value = 0.852224964839111
$ws.Range("M10").Value = 0.819389
$ws.Range("N10").Value = 2.458167
$ws.Range("O10").Value = 0.03318376186120772
$ws.Range("P10").Value = 0.03318376186120772
$ws.Range("Q10").Value = 0.3445751980030001
$ws.Range("R10").Value = 3.101176782027
$ws.Range("S10").Value = 0.02828003028539719
$ws.Range("T10").Value = 0.02828003028539719

# Row 11
$ws.Range("G11").Value = 0.420527
$ws.Range("H11").Value = 1.261581
$ws.Range("I11").Value = 0.852224964839111
$ws.Range("J11").Value = 0.852224964839111
$ws.Range("M11").Value = 3.080288333333333
$ws.Range("N11").Value = 9.240864999999999
$ws.Range("O11").Value = 0.1247460662971919
$ws.Range("P11").Value = 0.1247460662971919
$ws.Range("Q11").Value = 1.295344411951667
$ws.Range("R11").Value = 11.658099707565
$ws.Range("S11").Value = 0.1063117119639418
$ws.Range("T11").Value = 0.1063117119639418
